$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "2018-00069-BN-0"
$ws.Range("A3").Value = "2018-11170-BN-0"
$ws.Range("A4").Value = "2018-03371-BN-2"
$ws.Range("A5").Value = "2018-00572-BN-3"
$ws.Range("B2").Select()
